# Generate Report for Handoff
# Swap the old run UUID (3bdeb3e2-aa55-4695-ae1a-6f6a981800e6) for the new one
# (e4cc04c7-577e-4eb0-ad03-2278bc75be67) everywhere it is referenced, update the
# handoff-file checksums, and bump the "Latest ... Datetime" timestamps to match
# the newer handoff run.

$wb = $excel.ActiveWorkbook

$oldGuid = "3bdeb3e2-aa55-4695-ae1a-6f6a981800e6"
$newGuid = "e4cc04c7-577e-4eb0-ad03-2278bc75be67"

$oldHash = "df8725621100c3035f3ab9e861baf01da22714ce"
$newHash = "4bb023d865805d7209240a69a638340e7e3a8ed1"

$repoBase = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/df5e2c84ab193c1135607106d5a60cc105c6321b/e2e/"

# ---------------------------------------------------------------------------
# "Overview" sheet: File Name / Path And Name / Latest HO Xliff Generate Date
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = "$newGuid.md"

$newPathAndName = "e2e\" + $newGuid + ".md"
$wsOverview.Range("B2").Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), ($repoBase + $newGuid + ".md"), "", "", $newPathAndName)

$wsOverview.Range("G2").Value = "2016-08-23 17:02:38"

# ---------------------------------------------------------------------------
# "zh-cn" sheet: Source File Name / Latest Handoff File / Latest Handoff Datetime
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("A2").Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), ($repoBase + $newGuid + ".md"), "", "", "$newGuid.md")

$wsZhCn.Range("G2").Value = "$newGuid.$newHash.zh-cn.xlf"
$wsZhCn.Range("H2").Value = "2016-08-23 17:02:33"

# ---------------------------------------------------------------------------
# "de-de" sheet: Source File Name / Latest Handoff File / Latest Handback DateTime
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("A2").Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), ($repoBase + $newGuid + ".md"), "", "", "$newGuid.md")

$wsDeDe.Range("G2").Value = "$newGuid.$newHash.de-de.xlf"
$wsDeDe.Range("H2").Value = "2016-08-23 17:02:38"
